# Generate Report for Handoff
# - Removes the stale "27b4ad79-..." handoff row from every sheet (the
#   handoff for that file is gone; its hyperlinks/rows collapse and the
#   trailing ".localization-config" row shifts up into its place).
# - Marks the "1a61298b-..." file's status as "Ready for handoff" and
#   refreshes its Latest Handoff Datetime on the zh-cn/de-de sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Drop existing hyperlinks up front -- row deletion below does not retarget
# them, so they're rebuilt from scratch after the grid is correct again.
$ws.Cells.Hyperlinks.Delete()

# The "27b4ad79-...md" row is row 3; deleting it shifts the
# ".localization-config" row up from 4 to 3.
$ws.Rows.Item(3).Delete()

$ws.Cells.Item(2, 2).Value2 = "Ready for handoff"
$ws.Cells.Item(2, 3).Value2 = "Ready for handoff"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/08191acce60b636cd3283ada74ef54f794f14098/e2e/1a61298b-aa10-4043-b781-2015c4fbe842.md", "", "", "1a61298b-aa10-4043-b781-2015c4fbe842.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/08191acce60b636cd3283ada74ef54f794f14098/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Cells.Hyperlinks.Delete()
$ws.Rows.Item(3).Delete()

$ws.Cells.Item(2, 2).Value2 = "Ready for handoff"
$ws.Cells.Item(2, 4).Value2 = "2016-03-08 23:04:53"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/08191acce60b636cd3283ada74ef54f794f14098/e2e/1a61298b-aa10-4043-b781-2015c4fbe842.md", "", "", "1a61298b-aa10-4043-b781-2015c4fbe842.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0d3ce3d54c8a07a9351f2b544b4c440e46808297/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1a61298b-aa10-4043-b781-2015c4fbe842.3bb3d00a8adda4e1e285f5c028672018c13ecb08.zh-cn.xlf", "", "", "1a61298b-aa10-4043-b781-2015c4fbe842.3bb3d00a8adda4e1e285f5c028672018c13ecb08.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/324a37b17f2cd4a3004ebaa83c6544091394cb9e/e2e/1a61298b-aa10-4043-b781-2015c4fbe842.md", "", "", "1a61298b-aa10-4043-b781-2015c4fbe842.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/245c339a52592a0ad3ee3e8fa8dd18c5cab93652/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1a61298b-aa10-4043-b781-2015c4fbe842.3bb3d00a8adda4e1e285f5c028672018c13ecb08.zh-cn.xlf", "", "", "1a61298b-aa10-4043-b781-2015c4fbe842.3bb3d00a8adda4e1e285f5c028672018c13ecb08.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/08191acce60b636cd3283ada74ef54f794f14098/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Cells.Hyperlinks.Delete()
$ws.Rows.Item(3).Delete()

$ws.Cells.Item(2, 2).Value2 = "Ready for handoff"
$ws.Cells.Item(2, 4).Value2 = "2016-03-08 23:05:01"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/08191acce60b636cd3283ada74ef54f794f14098/e2e/1a61298b-aa10-4043-b781-2015c4fbe842.md", "", "", "1a61298b-aa10-4043-b781-2015c4fbe842.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b85dcb373392ceef6b614720f4cf7b7639bdf14a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1a61298b-aa10-4043-b781-2015c4fbe842.3bb3d00a8adda4e1e285f5c028672018c13ecb08.de-de.xlf", "", "", "1a61298b-aa10-4043-b781-2015c4fbe842.3bb3d00a8adda4e1e285f5c028672018c13ecb08.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4672d35c8bc12e8488654271ef4241f49cb190e7/e2e/1a61298b-aa10-4043-b781-2015c4fbe842.md", "", "", "1a61298b-aa10-4043-b781-2015c4fbe842.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0e679e6c4b0ba67bbb4dc9f9f02bbc0503437fcf/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1a61298b-aa10-4043-b781-2015c4fbe842.3bb3d00a8adda4e1e285f5c028672018c13ecb08.de-de.xlf", "", "", "1a61298b-aa10-4043-b781-2015c4fbe842.3bb3d00a8adda4e1e285f5c028672018c13ecb08.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/08191acce60b636cd3283ada74ef54f794f14098/.localization-config", "", "", ".localization-config")
